# Backlog update: add several new items and a "perfect client" entry,
# remove the now-unused "Agregar jugadores" / "Hoja3" sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("General")

# --- Row 7 changes: replaces "Poder agregar jugadores" level-logic item
#     with a new, shorter "Lógica de perfiles" item (no Description/Status). ---
$ws.Range("B7").Value = "Lógica de perfiles"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""

# --- Rows 8-19: same rows, content shifted up (the old rows 9-19 move to 8-18,
#     and a brand-new "perfect client" row becomes row 19). ---
$ws.Range("B8").Value = "Conexión de test para los clientes"
$ws.Range("C8").Value = "Hacer una url tipo hello world que tenga un help"
$ws.Range("D8").Value = "ok"

$ws.Range("B9").Value = "Definir arquitectura del server"
$ws.Range("C9").Value = "Crear modelo, crear tipos, crear main"
$ws.Range("D9").Value = "ok"

$ws.Range("B10").Value = "investigar sobre módulos en Node"
$ws.Range("C10").Value = "Queremos poner cada cosa como un módulo"
$ws.Range("D10").Value = "ok"

$ws.Range("B11").Value = "Cliente de prueba"
$ws.Range("C11").Value = "Que muestre lo que se le preguntó y conteste algún valor hardcodeado"
$ws.Range("D11").Value = "ok"

$ws.Range("B12").Value = "Hacer un log"
$ws.Range("C12").Value = "De lo que se le envía al servidor"
$ws.Range("D12").Value = "ok"

$ws.Range("B13").Value = "Despacahar las urls"
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = "ok"

$ws.Range("B14").Value = "Notificar a los jugadores"
$ws.Range("C14").Value = "A partir del player manager notificar a cada jugador con lo que corresponda"
$ws.Range("D14").Value = "ok"

$ws.Range("B15").Value = "Leader board"
$ws.Range("C15").Value = "Muestra los puntajes parciales de todos los jugadores"
$ws.Range("D15").Value = "ok"

$ws.Range("B16").Value = "Agregar los puntajes a los jugadores"
$ws.Range("C16").Value = "Cuando se agrega un jugador empieza con 0. Y se pueden sumar o restar"
$ws.Range("D16").Value = "ok"

$ws.Range("B17").Value = "Enviarle preguntas a los jugadores"
$ws.Range("C17").Value = "Enviarles un request con la pregunta a los jugadores"
$ws.Range("D17").Value = "ok"

$ws.Range("B18").Value = "Agregar más preguntas"
$ws.Range("C18").Value = "Y sus respuestas. También determinar el orden de los niveles."
$ws.Range("D18").Value = ""

$ws.Range("B19").Value = "Escribir un cliente perfecto"
$ws.Range("C19").Value = "Que conteste bien todas las preguntas y tenga puntaje perfecto."
$ws.Range("D19").Value = ""

# --- New backlog rows 20-29 ---
$ws.Range("A20").Value = 19
$ws.Range("B20").Value = "Ordenar el leaderboard"
$ws.Range("C20").Value = "Por puntaje"

$ws.Range("A21").Value = 20
$ws.Range("B21").Value = "Poner estado al juego"
$ws.Range("C21").Value = "No debe arrancar de una sino cuando se le da start"

$ws.Range("A22").Value = 21
$ws.Range("B22").Value = "Poner pausa"
$ws.Range("C22").Value = "Y resume"

$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "Hacer perfiles"
$ws.Range("C23").Value = "Además tener un perfil de práctica (sin puntos)"

$ws.Range("A24").Value = 23
$ws.Range("B24").Value = "Preguntas por query string"
$ws.Range("C24").Value = 'Por post de forma "pregunta=¿Qué día es hoy?"'

$ws.Range("A25").Value = 24
$ws.Range("B25").Value = "Agregar preguntas"
$ws.Range("C25").Value = "Colores primarios"

$ws.Range("A26").Value = 25
$ws.Range("B26").Value = "Cambiar el perfil mientras está en pausa"

$ws.Range("A27").Value = 26
$ws.Range("B27").Value = "Agregar botón para ir al leadearboard"
$ws.Range("C27").Value = "Después de haber sido agregado"

$ws.Range("A28").Value = 27
$ws.Range("B28").Value = 'Poder decir "paso"'
$ws.Range("C28").Value = "Y que reste puntos pero tira la siguiente pregunta"

$ws.Range("A29").Value = 28
$ws.Range("B29").Value = "Mejorar css"
$ws.Range("C29").Value = "A Adrián no le gustó"

# --- Column B grew a bit wider to fit the longest new entry ---
$ws.Columns("B").ColumnWidth = 37.333333333333336

# --- View: scroll down to the newly-added rows ---
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B30").Select() | Out-Null

# --- Drop the sheets that are no longer needed ---
$wb.Worksheets("Agregar jugadores").Delete() | Out-Null
$wb.Worksheets("Hoja3").Delete() | Out-Null
